$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.213901162147522
$ws.Range("B1").Value = 2.141748428344727
$ws.Range("C1").Value = 5.896011352539062
$ws.Range("D1").Value = 1.107722640037537
$ws.Range("E1").Value = 1.245447278022766
